$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to be treated as text so that values such as
# "1.00", "0.999" or "0.0000134" are not silently coerced into numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "57.649.81"
$ws.Range("E2").Value = "  -2.04%  "
$ws.Range("D3").Value = "2.445.72"
$ws.Range("E3").Value = "  -3.03%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "514.87"
$ws.Range("E5").Value = "  -4.13%  "
$ws.Range("E6").Value = "  -3.77%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "0.554"
$ws.Range("E8").Value = "  -2.25%  "
$ws.Range("D9").Value = "2.445.87"
$ws.Range("E9").Value = "  -3.02%  "
$ws.Range("D10").Value = "0.0976"
$ws.Range("E10").Value = "  -3.42%  "
$ws.Range("D11").Value = "0.156"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").Value = "5.23"
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").Value = "0.338"
$ws.Range("E13").Value = "  -2.75%  "
$ws.Range("D14").Value = "2.882.01"
$ws.Range("E14").Value = "  -2.75%  "
$ws.Range("D15").Value = "57.623.49"
$ws.Range("E15").Value = "  -2.11%  "
$ws.Range("D16").Value = "22.00"
$ws.Range("E16").Value = "  -4.25%  "
$ws.Range("D17").Value = "0.0000134"
$ws.Range("E17").Value = "  -3.14%  "
$ws.Range("D18").Value = "2.424.62"
$ws.Range("E18").Value = "  -3.63%  "
$ws.Range("D19").Value = "10.58"
$ws.Range("E19").Value = "  -4.60%  "
$ws.Range("D20").Value = "317.52"
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("D21").Value = "4.13"
$ws.Range("E21").Value = "  -2.93%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "5.68"
$ws.Range("E23").Value = "  -4.22%  "
$ws.Range("D24").Value = "63.99"
$ws.Range("E24").Value = "  -1.71%  "
$ws.Range("D25").Value = "0.404"
$ws.Range("E25").Value = "  -3.59%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  -3.23%  "
$ws.Range("D28").Value = "7.27"
$ws.Range("E28").Value = "  -3.02%  "
$ws.Range("D29").Value = "0.0₃0728"
$ws.Range("E29").Value = "  -4.85%  "
$ws.Range("D30").Value = "165.45"
$ws.Range("E30").Value = "  -3.33%  "
$ws.Range("E31").Value = "  -4.66%  "
$ws.Range("D32").Value = "6.17"
$ws.Range("E32").Value = "  -6.64%  "
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").Value = "17.94"
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("E37").Value = "  -7.55%  "
$ws.Range("E38").Value = "  -3.16%  "
$ws.Range("E39").Value = "  -5.07%  "
$ws.Range("D40").Value = "0.780"
$ws.Range("E40").Value = "  -3.54%  "
$ws.Range("D41").Value = "3.40"
$ws.Range("E41").Value = "  -4.98%  "
$ws.Range("D42").Value = "270.11"
$ws.Range("E42").Value = "  -4.89%  "
$ws.Range("D43").Value = "4.92"
$ws.Range("E43").Value = "  -4.80%  "
$ws.Range("D44").Value = "0.585"
$ws.Range("E44").Value = "  -3.60%  "
$ws.Range("D45").Value = "122.91"
$ws.Range("E45").Value = "  -5.83%  "
$ws.Range("D46").Value = "0.0903"
$ws.Range("E46").Value = "  -2.02%  "
$ws.Range("D47").Value = "0.0483"
$ws.Range("E47").Value = "  -4.21%  "
$ws.Range("D48").Value = "0.0208"
$ws.Range("E48").Value = "  -5.17%  "
$ws.Range("D49").Value = "16.57"
$ws.Range("E49").Value = "  -4.30%  "
$ws.Range("D50").Value = "1.714.77"
$ws.Range("E50").Value = "  -2.25%  "
$ws.Range("E51").Value = "  -2.46%  "

# Restore the default (no explicit style) look of the price column, matching
# the original workbook where these cells carry no style index.
$priceRange.Style = "Normal"
